$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.682.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "'1.922.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'240.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.96%  "

$ws.Range("D6").Value = "'0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "'0.4944"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.47%  "

$ws.Range("D8").Value = "'0.2996"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.22%  "

$ws.Range("D9").Value = "'0.06770"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("D10").Value = "'1.910.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.22%  "

$ws.Range("D11").Value = "'17.26"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").Value = "'0.07353"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("D13").Value = "'5.208"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.12%  "

$ws.Range("D14").Value = "'88.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.71%  "

$ws.Range("D15").Value = "'0.6760"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("D16").Value = "'30.657.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "

$ws.Range("D17").Value = "'0.000007961"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "'13.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.84%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").Value = "'2.171.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.87%  "

$ws.Range("D21").Value = "'5.418"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.21%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").Value = "'200.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.73%  "

$ws.Range("D24").Value = "'6.343"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.52%  "

$ws.Range("D25").Value = "'9.669"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.11%  "

$ws.Range("D26").Value = "'164.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.16%  "

$ws.Range("D27").Value = "'18.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.27%  "

$ws.Range("D28").Value = "'1.967"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.16%  "

$ws.Range("D29").Value = "'1.478"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.58%  "

$ws.Range("D30").Value = "'4.380"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").Value = "'0.09225"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.39%  "

$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").Value = "'0.05301"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.81%  "

$ws.Range("D34").Value = "'0.7457"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.18%  "

$ws.Range("D35").Value = "'1.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.96%  "

$ws.Range("D36").Value = "'2.716"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.06%  "

$ws.Range("D37").Value = "'0.01852"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "'2.724"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.43%  "

$ws.Range("D39").Value = "'0.9286"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("D40").Value = "'2.097"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").Value = "'0.4483"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("D42").Value = "'5.982"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.63%  "

$ws.Range("D43").Value = "'72.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +24.73%  "

$ws.Range("D44").Value = "'106.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("D46").Value = "'0.1403"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.13%  "

$ws.Range("D47").Value = "'7.673"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.75%  "

$ws.Range("D48").Value = "'9.037"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.49%  "

$ws.Range("D49").Value = "'35.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.53%  "

$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "'0.4050"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.71%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05884"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "

